$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 6: new task "Issue method to get from database" ---
$ws.Range("B6").Value = "Issue method to get from database"
$ws.Range("D6").Value = "Functional"
$ws.Range("F6").Value = "Ruwan"
$ws.Range("H6").Value = "1 day"
$ws.Range("J6").Value = "new"
$ws.Range("L6").Value = "Ruwan/Philip"

# --- Mark the two previous rows (4 & 5) as Done ---
$ws.Range("J4").Value = "Done"
$ws.Range("J5").Value = "Done"

# --- Row 7: new task "Add employees to service" ---
$ws.Range("B7").Value = "Add employees to service"
$ws.Range("D7").Value = "Functional"
$ws.Range("F7").Value = "Ruwan"
$ws.Range("H7").Value = "1 day"
$ws.Range("J7").Value = "new"
$ws.Range("L7").Value = "Ruwan"

# --- Row 8: new task "Remove service status from UI" ---
$ws.Range("D8").Value = "Modification"
$ws.Range("B8").Value = "Remove service status from UI"
$ws.Range("F8").Value = "Ruwan"
$ws.Range("H8").Value = ".5 h"
$ws.Range("J8").Value = "new"
$ws.Range("L8").Value = "Ruwan"

# --- Row 9: new task "Activate manager functions with authorized user name password" ---
$ws.Range("B9").Value = "Activate manager functions with authorized user name password"
$ws.Range("D9").Value = "Functional"
$ws.Range("F9").Value = "Ruwan"
$ws.Range("H9").Value = "1 day"
$ws.Range("J9").Value = "new"
$ws.Range("L9").Value = "Ruwan"

# --- Row 10: new task "Show reorder level reach with red color" ---
$ws.Range("B10").Value = "Show reorder level reach with red color"
$ws.Range("D10").Value = "Functional"
$ws.Range("F10").Value = "Ruwan"
$ws.Range("H10").Value = "1 day"
$ws.Range("J10").Value = "new "
$ws.Range("L10").Value = "Ruwan"

# --- Dates for the new rows (reuse the existing date style via copy/paste-special) ---
$ws.Range("N5").Copy()
$ws.Range("N6:N10").PasteSpecial(-4122)
$ws.Range("P5").Copy()
$ws.Range("P6:P10").PasteSpecial(-4122)

$ws.Range("N6").Value = 42981
$ws.Range("P6").Value = 42980
$ws.Range("N7").Value = 42981
$ws.Range("P7").Value = 42980
$ws.Range("N8").Value = 42981
$ws.Range("P8").Value = 42980
$ws.Range("N9").Value = 42983
$ws.Range("P9").Value = 42980
$ws.Range("N10").Value = 42983
$ws.Range("P10").Value = 42980

# --- Wrap the long text in B9 and grow row 9 to fit two lines ---
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 30

# --- Column widths: B (task name) wider, J (status) wider to fit "Done" ---
$ws.Columns.Item(2).ColumnWidth = 35.83
$ws.Columns.Item(10).ColumnWidth = 11.5

# --- Selection ends on the last edited cell ---
$ws.Range("P10").Select() | Out-Null
